$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should match the existing
# header formatting (bold font + border, style used by B1:H1). Copy the
# format from H1 (the last existing header cell) over to I1:J1, then set
# the text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-16.
$data = @(
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 9),
    @(6, 7),
    @(10, 10),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(2, 2),
    @(9, 9),
    @(3, 3),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
